# Auto-generated Excel COM-interop script
# Applies the Famfrit_Profits price-refresh values from the scheduled runner diff.
$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
# Row 17
$ws.Range("H17").Value = 742759.9
$ws.Range("J17").Value = 742759.9
$ws.Range("L17").Value = 2228279.7
$ws.Range("N17").Value = -2228615.7
# Row 113
$ws.Range("H113").Value = 6254.25
$ws.Range("I113").Value = 4553.909
$ws.Range("J113").Value = 8332.444
$ws.Range("K113").Value = 4553.909
$ws.Range("L113").Value = 8332.444
$ws.Range("M113").Value = -1299.909
$ws.Range("N113").Value = -14840.444
# Row 116
$ws.Range("H116").Value = 5183.769
$ws.Range("I116").Value = 4798.625
$ws.Range("K116").Value = 4798.625
$ws.Range("M116").Value = -1356.625
# Row 131
$ws.Range("H131").Value = 5219.222
$ws.Range("I131").Value = 4105.1113
$ws.Range("K131").Value = 12315.3339
$ws.Range("M131").Value = -7275.333899999998
# Row 132
$ws.Range("H132").Value = 4372.278
$ws.Range("I132").Value = 4466
$ws.Range("K132").Value = 13398
$ws.Range("M132").Value = -10868

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
# Row 2
$ws.Range("H2").Value = 736.3
$ws.Range("I2").Value = 600.7778
$ws.Range("J2").Value = 1956
$ws.Range("K2").Value = 600.7778
$ws.Range("L2").Value = 1956
$ws.Range("M2").Value = -487.7778
$ws.Range("N2").Value = -2182
# Row 32
$ws.Range("H32").Value = 13703069
$ws.Range("I32").Value = 15627697
$ws.Range("K32").Value = 15627697
$ws.Range("M32").Value = -15627410
# Row 102
$ws.Range("H102").Value = 2350.818
$ws.Range("J102").Value = 5394.6
$ws.Range("L102").Value = 5394.6
$ws.Range("N102").Value = -8638.6
# Row 110
$ws.Range("H110").Value = 17120.666
$ws.Range("I110").Value = 18781.238
$ws.Range("K110").Value = 18781.238
$ws.Range("M110").Value = -16736.238
# Row 116
$ws.Range("H116").Value = 736.3
$ws.Range("I116").Value = 600.7778
$ws.Range("J116").Value = 1956
$ws.Range("K116").Value = 600.7778
$ws.Range("L116").Value = 1956
$ws.Range("M116").Value = 1693.2222
$ws.Range("N116").Value = -6544

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
# Row 3
$ws.Range("H3").Value = 736.3
$ws.Range("I3").Value = 600.7778
$ws.Range("J3").Value = 1956
$ws.Range("K3").Value = 600.7778
$ws.Range("L3").Value = 1956
$ws.Range("M3").Value = -486.7778
$ws.Range("N3").Value = -2184
# Row 86
$ws.Range("H86").Value = 16754.047
$ws.Range("I86").Value = 8646.75
$ws.Range("K86").Value = 8646.75
$ws.Range("M86").Value = -7523.75
# Row 89
$ws.Range("H89").Value = 16754.047
$ws.Range("I89").Value = 8646.75
$ws.Range("K89").Value = 43233.75
$ws.Range("M89").Value = -37617.75

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
# Row 68
$ws.Range("H68").Value = 40000
$ws.Range("J68").Value = 40000
$ws.Range("L68").Value = 40000
$ws.Range("N68").Value = -41498
# Row 71
$ws.Range("H71").Value = 40000
$ws.Range("J71").Value = 40000
$ws.Range("L71").Value = 120000
$ws.Range("N71").Value = -127488
# Row 94
$ws.Range("H94").Value = 1464.4117
$ws.Range("J94").Value = 1511.6364
$ws.Range("L94").Value = 1511.6364
$ws.Range("N94").Value = -2413.6364
# Row 99
$ws.Range("H99").Value = 15273.125
$ws.Range("I99").Value = 15273.125
$ws.Range("J99").Value = 0
$ws.Range("K99").Value = 15273.125
$ws.Range("L99").Value = 0
$ws.Range("M99").Value = -13775.125
$ws.Range("N99").ClearContents()
# Row 105
$ws.Range("H105").Value = 7296.85
$ws.Range("I105").Value = 1677.9286
$ws.Range("J105").Value = 20407.666
$ws.Range("K105").Value = 1677.9286
$ws.Range("L105").Value = 20407.666
$ws.Range("M105").Value = 69.07140000000004
$ws.Range("N105").Value = -23901.666
# Row 126
$ws.Range("H126").Value = 15273.125
$ws.Range("I126").Value = 15273.125
$ws.Range("J126").Value = 0
$ws.Range("K126").Value = 45819.375
$ws.Range("L126").Value = 0
$ws.Range("M126").Value = -43349.375
$ws.Range("N126").ClearContents()

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
# Row 7
$ws.Range("H7").Value = 683.5
$ws.Range("I7").Value = 740.2
$ws.Range("K7").Value = 2220.6
$ws.Range("M7").Value = -2108.6
# Row 18
$ws.Range("H18").Value = 1764.5834
$ws.Range("I18").Value = 808.3333
$ws.Range("J18").Value = 4633.3335
$ws.Range("K18").Value = 2424.9999
$ws.Range("L18").Value = 13900.0005
$ws.Range("M18").Value = -2255.9999
$ws.Range("N18").Value = -14238.0005
# Row 19
$ws.Range("H19").Value = 1354.8889
$ws.Range("I19").Value = 313.42856
$ws.Range("J19").Value = 5000
$ws.Range("K19").Value = 940.28568
$ws.Range("L19").Value = 15000
$ws.Range("M19").Value = -766.28568
$ws.Range("N19").Value = -15348
# Row 23
$ws.Range("H23").Value = 1360.2
$ws.Range("J23").Value = 1566.6666
$ws.Range("L23").Value = 4699.9998
$ws.Range("N23").Value = -5169.9998
# Row 25
$ws.Range("H25").Value = 1833.3334
$ws.Range("I25").Value = 200
$ws.Range("K25").Value = 600
$ws.Range("M25").Value = -431
# Row 30
$ws.Range("H30").Value = 1833.3334
$ws.Range("I30").Value = 200
$ws.Range("K30").Value = 600
$ws.Range("M30").Value = -498
# Row 56
$ws.Range("H56").Value = 23758.3
$ws.Range("I56").Value = 23758.3
$ws.Range("K56").Value = 23758.3
$ws.Range("M56").Value = -23228.3
# Row 76
$ws.Range("H76").Value = 3014
$ws.Range("I76").Value = 2013
$ws.Range("J76").Value = 4015
$ws.Range("K76").Value = 6039
$ws.Range("L76").Value = 12045
$ws.Range("M76").Value = -5656
$ws.Range("N76").Value = -12811
# Row 79
$ws.Range("H79").Value = 3014
$ws.Range("I79").Value = 2013
$ws.Range("J79").Value = 4015
$ws.Range("K79").Value = 6039
$ws.Range("L79").Value = 12045
$ws.Range("M79").Value = -4713
$ws.Range("N79").Value = -14697
# Row 140
$ws.Range("H140").Value = 2285.5
$ws.Range("I140").Value = 1464.8334
$ws.Range("K140").Value = 4394.5002
$ws.Range("M140").Value = 785.4997999999996

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
# Row 32
$ws.Range("H32").Value = 44000
$ws.Range("I32").Value = 44000
$ws.Range("K32").Value = 44000
$ws.Range("M32").Value = -43704
# Row 42
$ws.Range("H42").Value = 29000.334
$ws.Range("J42").Value = 27000
$ws.Range("L42").Value = 27000
$ws.Range("N42").Value = -27970
# Row 86
$ws.Range("H86").Value = 27400
$ws.Range("J86").Value = 27400
$ws.Range("L86").Value = 27400
$ws.Range("N86").Value = -29772
# Row 89
$ws.Range("H89").Value = 27400
$ws.Range("J89").Value = 27400
$ws.Range("L89").Value = 82200
$ws.Range("N89").Value = -94056
# Row 115
$ws.Range("H115").Value = 29000.334
$ws.Range("J115").Value = 27000
$ws.Range("L115").Value = 27000
$ws.Range("N115").Value = -29350
# Row 122
$ws.Range("H122").Value = 7107.7085
$ws.Range("I122").Value = 3088.5715
$ws.Range("K122").Value = 9265.7145
$ws.Range("M122").Value = -6815.7145
# Row 136
$ws.Range("H136").Value = 25148.268
$ws.Range("J136").Value = 25148.268
$ws.Range("L136").Value = 75444.804
$ws.Range("N136").Value = -80544.804

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
# Row 55
$ws.Range("H55").Value = 484.8889
$ws.Range("I55").Value = 228.25
$ws.Range("K55").Value = 228.25
$ws.Range("M55").Value = -55.25
# Row 61
$ws.Range("H61").Value = 4756.2856
$ws.Range("I61").Value = 3941.1428
$ws.Range("J61").Value = 5571.4287
$ws.Range("K61").Value = 3941.1428
$ws.Range("L61").Value = 5571.4287
$ws.Range("M61").Value = -3739.1428
$ws.Range("N61").Value = -5975.4287
# Row 113
$ws.Range("H113").Value = 4756.2856
$ws.Range("I113").Value = 3941.1428
$ws.Range("J113").Value = 5571.4287
$ws.Range("K113").Value = 3941.1428
$ws.Range("L113").Value = 5571.4287
$ws.Range("M113").Value = -1771.1428
$ws.Range("N113").Value = -9911.4287

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
# Row 27
$ws.Range("H27").Value = 74499.5
$ws.Range("J27").Value = 74499.5
$ws.Range("L27").Value = 74499.5
$ws.Range("N27").Value = -74637.5
# Row 115
$ws.Range("H115").Value = 60749.25
$ws.Range("I115").Value = 50000
$ws.Range("J115").Value = 64332.332
$ws.Range("K115").Value = 50000
$ws.Range("L115").Value = 64332.332
$ws.Range("M115").Value = -48433
$ws.Range("N115").Value = -67466.33199999999
